# Scheduled-runner refresh: update FFXIV market-board price columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
# across the eight crafting-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 207.5
$ws.Range("I9").Value = 148
$ws.Range("J9").Value = 267
$ws.Range("K9").Value = 148
$ws.Range("L9").Value = 267
$ws.Range("M9").Value = 21
$ws.Range("N9").Value = -605
$ws.Range("H33").Value = 144.09091
$ws.Range("I33").Value = 170.5
$ws.Range("J33").Value = 73.666664
$ws.Range("K33").Value = 170.5
$ws.Range("L33").Value = 73.666664
$ws.Range("M33").Value = 58.5
$ws.Range("N33").Value = -531.666664
$ws.Range("H55").Value = 343.375
$ws.Range("J55").Value = 537.75
$ws.Range("L55").Value = 537.75
$ws.Range("N55").Value = -965.75
$ws.Range("H74").Value = 4500
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1564
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4500
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -7820
$ws.Range("N77").Value = -34360
$ws.Range("H132").Value = 114303.875
$ws.Range("I132").Value = 283076.03
$ws.Range("J132").Value = 14290.741
$ws.Range("K132").Value = 849228.0900000001
$ws.Range("L132").Value = 42872.223
$ws.Range("M132").Value = -846698.0900000001
$ws.Range("N132").Value = -47932.223
$ws.Range("H137").Value = 4231.6
$ws.Range("I137").Value = 1799.25
$ws.Range("J137").Value = 5116.091
$ws.Range("K137").Value = 5397.75
$ws.Range("L137").Value = 15348.273
$ws.Range("M137").Value = -2847.75
$ws.Range("N137").Value = -20448.273
$ws.Range("H138").Value = 5522.965
$ws.Range("J138").Value = 7082.073
$ws.Range("L138").Value = 21246.219
$ws.Range("N138").Value = -31526.219
$ws.Range("H141").Value = 7359
$ws.Range("I141").Value = 6183.846
$ws.Range("K141").Value = 18551.538
$ws.Range("M141").Value = -13371.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 858572.3
$ws.Range("I2").Value = 1589925.9
$ws.Range("K2").Value = 1589925.9
$ws.Range("M2").Value = -1589812.9
$ws.Range("H45").Value = 2237.923
$ws.Range("I45").Value = 2239.25
$ws.Range("K45").Value = 2239.25
$ws.Range("M45").Value = -1862.25
$ws.Range("H61").Value = 21325.125
$ws.Range("I61").Value = 23791.6
$ws.Range("K61").Value = 23791.6
$ws.Range("M61").Value = -23579.6
$ws.Range("H116").Value = 858572.3
$ws.Range("I116").Value = 1589925.9
$ws.Range("K116").Value = 1589925.9
$ws.Range("M116").Value = -1587631.9
$ws.Range("H132").Value = 27288.691
$ws.Range("I132").Value = 34612.234
$ws.Range("J132").Value = 13455.333
$ws.Range("K132").Value = 103836.702
$ws.Range("L132").Value = 40365.999
$ws.Range("M132").Value = -101306.702
$ws.Range("N132").Value = -45425.999
$ws.Range("H136").Value = 21325.125
$ws.Range("I136").Value = 23791.6
$ws.Range("K136").Value = 71374.79999999999
$ws.Range("M136").Value = -68824.79999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 858572.3
$ws.Range("I3").Value = 1589925.9
$ws.Range("K3").Value = 1589925.9
$ws.Range("M3").Value = -1589811.9
$ws.Range("H20").Value = 3781.3333
$ws.Range("I20").Value = 3520.5
$ws.Range("K20").Value = 3520.5
$ws.Range("M20").Value = -3273.5
$ws.Range("H94").Value = 457288.94
$ws.Range("I94").Value = 1054046.4
$ws.Range("K94").Value = 1054046.4
$ws.Range("M94").Value = -1053595.4
$ws.Range("H105").Value = 2390.3333
$ws.Range("I105").Value = 2472.3635
$ws.Range("J105").Value = 2333.9375
$ws.Range("K105").Value = 2472.3635
$ws.Range("L105").Value = 2333.9375
$ws.Range("M105").Value = -725.3634999999999
$ws.Range("N105").Value = -5827.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 3600
$ws.Range("I10").Value = 2533.3333
$ws.Range("J10").Value = 4666.6665
$ws.Range("K10").Value = 2533.3333
$ws.Range("L10").Value = 4666.6665
$ws.Range("M10").Value = -2394.3333
$ws.Range("N10").Value = -4944.6665
$ws.Range("H11").Value = 2160.6
$ws.Range("J11").Value = 3533.3333
$ws.Range("L11").Value = 3533.3333
$ws.Range("N11").Value = -3813.3333
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = $null
$ws.Range("H15").Value = 6000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = $null
$ws.Range("H31").Value = 2776
$ws.Range("I31").Value = 1030.6666
$ws.Range("K31").Value = 1030.6666
$ws.Range("M31").Value = -735.6666
$ws.Range("H34").Value = 2776
$ws.Range("I34").Value = 1030.6666
$ws.Range("K34").Value = 1030.6666
$ws.Range("M34").Value = -828.6666
$ws.Range("H50").Value = 92000
$ws.Range("J50").Value = 92000
$ws.Range("L50").Value = 92000
$ws.Range("N50").Value = -93250
$ws.Range("H62").Value = 63141.2
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 63141.2
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("H99").Value = 18411.777
$ws.Range("J99").Value = 8499.833000000001
$ws.Range("L99").Value = 8499.833000000001
$ws.Range("N99").Value = -11495.833
$ws.Range("H122").Value = 8120.875
$ws.Range("I122").Value = 7989
$ws.Range("K122").Value = 23967
$ws.Range("M122").Value = -21517
$ws.Range("H126").Value = 18411.777
$ws.Range("J126").Value = 8499.833000000001
$ws.Range("L126").Value = 25499.499
$ws.Range("N126").Value = -30439.499
$ws.Range("H132").Value = 8339116
$ws.Range("I132").Value = 9262458
$ws.Range("K132").Value = 27787374
$ws.Range("M132").Value = -27784844
$ws.Range("H134").Value = 2176.4048
$ws.Range("I134").Value = 2130.225
$ws.Range("K134").Value = 6390.674999999999
$ws.Range("M134").Value = -3855.674999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 230
$ws.Range("I38").Value = 118.75
$ws.Range("K38").Value = 356.25
$ws.Range("M38").Value = -9.25
$ws.Range("H56").Value = 6670.864
$ws.Range("I56").Value = 6670.864
$ws.Range("K56").Value = 6670.864
$ws.Range("M56").Value = -6140.864
$ws.Range("H97").Value = 220
$ws.Range("I97").Value = 402
$ws.Range("J97").Value = 174.5
$ws.Range("K97").Value = 1206
$ws.Range("L97").Value = 523.5
$ws.Range("M97").Value = -710
$ws.Range("N97").Value = -1515.5
$ws.Range("H132").Value = 2547.375
$ws.Range("I132").Value = 1482.7142
$ws.Range("K132").Value = 13344.4278
$ws.Range("M132").Value = -10814.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7693.364
$ws.Range("I102").Value = 7329.1577
$ws.Range("K102").Value = 7329.1577
$ws.Range("M102").Value = -5707.1577
$ws.Range("H107").Value = 15873772
$ws.Range("J107").Value = 1134.5
$ws.Range("L107").Value = 1134.5
$ws.Range("N107").Value = -4974.5
$ws.Range("H122").Value = 735152
$ws.Range("J122").Value = 2015.4
$ws.Range("L122").Value = 6046.200000000001
$ws.Range("N122").Value = -10946.2
$ws.Range("H132").Value = 5243.125
$ws.Range("J132").Value = 5271.3335
$ws.Range("L132").Value = 15814.0005
$ws.Range("N132").Value = -20874.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5130.531
$ws.Range("I7").Value = 3807.5293
$ws.Range("K7").Value = 3807.5293
$ws.Range("M7").Value = -3695.5293
$ws.Range("H16").Value = 7693963
$ws.Range("I16").Value = 10001182
$ws.Range("J16").Value = 3233.5
$ws.Range("K16").Value = 10001182
$ws.Range("L16").Value = 3233.5
$ws.Range("M16").Value = -10001012
$ws.Range("N16").Value = -3573.5
$ws.Range("H40").Value = 5338.3076
$ws.Range("I40").Value = 4387.4443
$ws.Range("K40").Value = 4387.4443
$ws.Range("M40").Value = -4251.4443
$ws.Range("H55").Value = 222.5
$ws.Range("J55").Value = 211.66667
$ws.Range("L55").Value = 211.66667
$ws.Range("N55").Value = -557.6666700000001
$ws.Range("H68").Value = 990870.0600000001
$ws.Range("I68").Value = 1750143.2
$ws.Range("K68").Value = 1750143.2
$ws.Range("M68").Value = -1749394.2
$ws.Range("H71").Value = 990870.0600000001
$ws.Range("I71").Value = 1750143.2
$ws.Range("K71").Value = 8750716
$ws.Range("M71").Value = -8746972
$ws.Range("H93").Value = 752.6087
$ws.Range("I93").Value = 714.125
$ws.Range("J93").Value = 840.5714
$ws.Range("K93").Value = 714.125
$ws.Range("L93").Value = 840.5714
$ws.Range("M93").Value = 533.875
$ws.Range("N93").Value = -3336.5714
$ws.Range("H122").Value = 54427250
$ws.Range("J122").Value = 15878245
$ws.Range("L122").Value = 47634735
$ws.Range("N122").Value = -47639635
$ws.Range("H126").Value = 5130.531
$ws.Range("I126").Value = 3807.5293
$ws.Range("K126").Value = 11422.5879
$ws.Range("M126").Value = -8952.5879
$ws.Range("H132").Value = 4052.3833
$ws.Range("I132").Value = 2960.9534
$ws.Range("K132").Value = 8882.860199999999
$ws.Range("M132").Value = -6352.860199999999
$ws.Range("H136").Value = 4162
$ws.Range("I136").Value = 4116.1665
$ws.Range("J136").Value = 4299.5
$ws.Range("K136").Value = 12348.4995
$ws.Range("L136").Value = 12898.5
$ws.Range("M136").Value = -9798.499500000002
$ws.Range("N136").Value = -17998.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 37714.5
$ws.Range("J46").Value = 37714.5
$ws.Range("L46").Value = 37714.5
$ws.Range("N46").Value = -38176.5
$ws.Range("H96").Value = 5319.778
$ws.Range("I96").Value = 3834
$ws.Range("K96").Value = 3834
$ws.Range("M96").Value = -2461
$ws.Range("H122").Value = 3001.8542
$ws.Range("I122").Value = 2186.8975
$ws.Range("K122").Value = 6560.6925
$ws.Range("M122").Value = -4110.6925
$ws.Range("H126").Value = 2016.2
$ws.Range("I126").Value = 1895.25
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 5685.75
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -3215.75
$ws.Range("N126").Value = -12440
$ws.Range("H134").Value = 37714.5
$ws.Range("J134").Value = 37714.5
$ws.Range("L134").Value = 113143.5
$ws.Range("N134").Value = -118213.5
$ws.Range("H136").Value = 7806.0645
$ws.Range("I136").Value = 3643.524
$ws.Range("J136").Value = 9020.138999999999
$ws.Range("K136").Value = 10930.572
$ws.Range("L136").Value = 27060.417
$ws.Range("M136").Value = -8380.572
$ws.Range("N136").Value = -32160.417
